# Auto-generated edit script: applies cryptocurrency price/volume/coin updates
# as described in the commit "Updated symbol list on Mon Feb 13 16:59:52 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ws, $ref, $val) {
    $rng = $ws.Range($ref)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextCell $ws "D2" "285.46"
Set-TextCell $ws "E2" "-10.48%"
Set-TextCell $ws "D3" "39.92"
Set-TextCell $ws "E3" "-3.32%"
Set-TextCell $ws "D4" "5.030"
Set-TextCell $ws "E4" "-3.98%"
Set-TextCell $ws "D5" "0.07280"
Set-TextCell $ws "E5" "-5.91%"
Set-TextCell $ws "D6" "4.300"
Set-TextCell $ws "E6" "-0.37%"
Set-TextCell $ws "D7" "1.518"
Set-TextCell $ws "E7" "-11.14%"
Set-TextCell $ws "D8" "0.9174"
Set-TextCell $ws "E8" "-3.66%"
Set-TextCell $ws "D9" "0.1196"
Set-TextCell $ws "E9" "-5.46%"
Set-TextCell $ws "D10" "0.1710"
Set-TextCell $ws "E10" "-6.75%"
Set-TextCell $ws "D11" "0.08646"
Set-TextCell $ws "E11" "-5.63%"
Set-TextCell $ws "D12" "0.04168"
Set-TextCell $ws "E12" "-4.51%"
Set-TextCell $ws "D13" "0.1051"
Set-TextCell $ws "E13" "-0.14%"
Set-TextCell $ws "D14" "0.001268"
Set-TextCell $ws "E14" "-0.98%"
Set-TextCell $ws "D15" "0.005964"
Set-TextCell $ws "E15" "-0.90%"
Set-TextCell $ws "B16" "LEO"
Set-TextCell $ws "C16" "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextCell $ws "D16" "3.400"
Set-TextCell $ws "E16" "1.33%"
Set-TextCell $ws "B17" "BTSEToken"
Set-TextCell $ws "C17" "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextCell $ws "D17" "2.397"
Set-TextCell $ws "E17" "-1.16%"
Set-TextCell $ws "B18" "BitpandaEcosystemToken"
Set-TextCell $ws "C18" "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
Set-TextCell $ws "D18" "0.3281"
Set-TextCell $ws "E18" "-2.14%"
Set-TextCell $ws "B19" "MCDex"
Set-TextCell $ws "C19" "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextCell $ws "D19" "7.846"
Set-TextCell $ws "E19" "1.90%"
Set-TextCell $ws "B20" "ProBitToken"
Set-TextCell $ws "C20" "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
Set-TextCell $ws "D20" "0.1342"
Set-TextCell $ws "E20" "-0.68%"
Set-TextCell $ws "B21" "ZBToken"
Set-TextCell $ws "C21" "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
Set-TextCell $ws "D21" "0.2887"
Set-TextCell $ws "E21" "2.43%"
Set-TextCell $ws "B22" "CoinExToken"
Set-TextCell $ws "C22" "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextCell $ws "D22" "0.03849"
Set-TextCell $ws "E22" "-4.50%"
Set-TextCell $ws "B23" "BitKan"
Set-TextCell $ws "C23" "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
Set-TextCell $ws "D23" "0.001272"
Set-TextCell $ws "E23" "0.49%"
Set-TextCell $ws "B24" "HotbitToken"
Set-TextCell $ws "C24" "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
Set-TextCell $ws "D24" "0.003790"
Set-TextCell $ws "E24" "-8.00%"
Set-TextCell $ws "D25" "0.0001284"
Set-TextCell $ws "E25" "1.22%"
Set-TextCell $ws "D26" "0.0003728"
Set-TextCell $ws "D38" "0.02318"
Set-TextCell $ws "E38" "-9.02%"
Set-TextCell $ws "D39" "0.04978"
Set-TextCell $ws "E39" "-7.07%"
Set-TextCell $ws "D40" "0.006224"
Set-TextCell $ws "E40" "216.12%"
Set-TextCell $ws "D41" "0.007704"
Set-TextCell $ws "E41" "-1.02%"
Set-TextCell $ws "D42" "0.1267"
Set-TextCell $ws "E42" "-4.04%"
Set-TextCell $ws "D43" "0.007372"
Set-TextCell $ws "E43" "0.23%"
Set-TextCell $ws "D44" "0.007435"
Set-TextCell $ws "E44" "-1.88%"
Set-TextCell $ws "D45" "0.3081"
Set-TextCell $ws "E45" "-10.48%"
Set-TextCell $ws "D46" "0.00006448"
Set-TextCell $ws "E46" "-3.99%"
Set-TextCell $ws "E47" "0.43%"
Set-TextCell $ws "E48" "15.02%"
Set-TextCell $ws "E49" "0.17%"
Set-TextCell $ws "D50" "0.00002106"
Set-TextCell $ws "E50" "0.43%"
Set-TextCell $ws "D51" "0.0002006"
Set-TextCell $ws "E51" "0.43%"
